$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("C16").Value = "45466414"
    $ws.Range("D16").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E16").Value = "2204"
    $ws.Range("F16").Value = 40000
    $ws.Range("C17").Value = "45466414"
    $ws.Range("D17").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E17").Value = "2205"
    $ws.Range("F17").Value = 40000
    $ws.Range("C18").Value = "45466414"
    $ws.Range("D18").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E18").Value = "2206"
    $ws.Range("F18").Value = 40000
    $ws.Range("C19").Value = "45466414"
    $ws.Range("D19").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E19").Value = "2207"
    $ws.Range("F19").Value = 40000
    $ws.Range("C20").Value = "45466414"
    $ws.Range("D20").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E20").Value = "2208"
    $ws.Range("F20").Value = 40000
    $ws.Range("C21").Value = "1041971915"
    $ws.Range("D21").Value = "ISABELLA BARROSO CANTILLO"
    $ws.Range("E21").Value = "2208"
    $ws.Range("F21").Value = 40000
    $ws.Range("C22").Value = "45466414"
    $ws.Range("D22").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E22").Value = "2209"
    $ws.Range("F22").Value = 40000
    $ws.Range("C23").Value = "1041971915"
    $ws.Range("D23").Value = "ISABELLA BARROSO CANTILLO"
    $ws.Range("E23").Value = "2209"
    $ws.Range("F23").Value = 40000
    $ws.Range("C24").Value = "45466414"
    $ws.Range("D24").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E24").Value = "2210"
    $ws.Range("F24").Value = 40000
    $ws.Range("C25").Value = "1041971915"
    $ws.Range("D25").Value = "ISABELLA BARROSO CANTILLO"
    $ws.Range("E25").Value = "2210"
    $ws.Range("F25").Value = 40000
    $ws.Range("C26").Value = "45466414"
    $ws.Range("D26").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E26").Value = "2211"
    $ws.Range("F26").Value = 40000
    $ws.Range("C27").Value = "1041971915"
    $ws.Range("D27").Value = "ISABELLA BARROSO CANTILLO"
    $ws.Range("E27").Value = "2211"
    $ws.Range("F27").Value = 40000
    $ws.Range("C28").Value = "45466414"
    $ws.Range("D28").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E28").Value = "2212"
    $ws.Range("F28").Value = 40000
    $ws.Range("C29").Value = "1041971915"
    $ws.Range("D29").Value = "ISABELLA BARROSO CANTILLO"
    $ws.Range("E29").Value = "2212"
    $ws.Range("F29").Value = 40000
    $ws.Range("C30").Value = "45466414"
    $ws.Range("D30").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E30").Value = "2301"
    $ws.Range("F30").Value = 40000
    $ws.Range("C31").Value = "1041971915"
    $ws.Range("D31").Value = "ISABELLA BARROSO CANTILLO"
    $ws.Range("E31").Value = "2301"
    $ws.Range("F31").Value = 40000
    $ws.Range("C32").Value = "45466414"
    $ws.Range("D32").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E32").Value = "2302"
    $ws.Range("F32").Value = 40000
    $ws.Range("C33").Value = "45466414"
    $ws.Range("D33").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E33").Value = "2303"
    $ws.Range("F33").Value = 40000
    $ws.Range("C34").Value = "45466414"
    $ws.Range("D34").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E34").Value = "2304"
    $ws.Range("F34").Value = 40000
    $ws.Range("C35").Value = "45466414"
    $ws.Range("D35").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E35").Value = "2305"
    $ws.Range("F35").Value = 40000
    $ws.Range("C36").Value = "45466414"
    $ws.Range("D36").Value = "IGNACIA MARIA FERNANDEZ CAICEDO"
    $ws.Range("E36").Value = "2306"
    $ws.Range("F36").Value = 29333